# Refresh the crypto price/volume snapshot (cryptos list update, GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "25.974.13" or "0.2630" that must stay text
# (they are not valid numbers, or have significant trailing zeros). A leading apostrophe
# forces Excel to store the literal text instead of re-typing it as a Number.

$ws.Range("D2").Value = '25.974.13'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '1.663.50'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''215.51'
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").Value = '''0.5074'
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''0.2630'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").Value = '''0.06372'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").Value = '''21.80'
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").Value = '''0.07394'
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").Value = '1.666.94'
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("D13").Value = '''4.485'
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '''0.5798'
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("D15").Value = '''0.000008440'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '''64.00'
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("D17").Value = '26.045.04'
$ws.Range("E17").Value = '  -2.35%  '
$ws.Range("D18").Value = '''4.897'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '''10.63'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").Value = '''188.58'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").Value = '''6.185'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '''145.10'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''7.560'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").Value = '''0.1186'
$ws.Range("E26").Value = '  +2.86%  '
$ws.Range("D27").Value = '''15.57'
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").Value = '''0.06546'
$ws.Range("E28").Value = '  +14.92%  '
$ws.Range("D29").Value = '''1.307'
$ws.Range("E29").Value = '  -1.40%  '
$ws.Range("D30").Value = '''1.313'
$ws.Range("E30").Value = '  -1.94%  '
$ws.Range("D31").Value = '''3.517'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '''3.491'
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("D33").Value = '''1.621'
$ws.Range("E33").Value = '  -2.71%  '
$ws.Range("D34").Value = '''1.016'
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("D35").Value = '''0.6043'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '''2.369'
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = '''2.682'
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").Value = '''6.199'
$ws.Range("E38").Value = '  +5.26%  '
$ws.Range("D39").Value = '''0.01604'
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("D40").Value = '1.074.00'
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").Value = '''0.8580'
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("E44").Value = '  +8.04%  '
$ws.Range("D45").Value = '1.811.52'
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").Value = '''56.00'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").Value = '''1.003'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").Value = '''7.988'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").Value = '''0.05205'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").Value = '''0.4299'
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("D51").Value = '''5.926'
$ws.Range("E51").Value = '  +2.53%  '
